$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Row 4 ("Maven et Eclipse"): answer text (E4) gains a leading
#    newline, and the row height becomes 127.5 (auto-fit, not an
#    explicit/custom height anymore).
# ------------------------------------------------------------------
$ws.Range("E4").Value = "`n<plugin>`n   <groupId>org.apache.maven.plugins</groupId>`n   <artifactId>maven-compiler-plugin</artifactId>`n   <configuration>`n      <source>1.8</source>`n      <target>1.8</target>`n   </configuration>`n</plugin> `n`n`n"
$ws.Rows(4).RowHeight = 127.5

# ------------------------------------------------------------------
# 2) Insert a brand-new row 5 ("Choix du HTML") right after row 4.
#    Copy the formatting of an existing plain/blank data row (row 13)
#    so the new cells land on the usual styles (s=5,6,5,5). The
#    answer cell (D5) is filled in further below, after the other two
#    new question rows, to mirror how the sheet was actually edited.
# ------------------------------------------------------------------
$ws.Rows(5).Insert()
$ws.Range("B13:E13").Copy()
$ws.Range("B5:E5").PasteSpecial(-4122)
$ws.Range("B5").Value = "Choix du HTML"
$ws.Range("C5").Value = "Faut-il choisir HTML5 ou XHTML1.0 Transitional dans les jsp ?"
$ws.Rows(5).RowHeight = 33

# ------------------------------------------------------------------
# 3) Rows 11 and 12 (currently two of the blank filler rows, already
#    styled s=5,6,5,5) become two new question rows: "Log4j2.xml" and
#    "Test d'un controller Servlet".
# ------------------------------------------------------------------
$ws.Range("B11").Value = "Log4j2.xml"
$ws.Range("C11").Value = "A quoi sert et comment utiliser le Log4j2.xml ?"
$ws.Range("D11").Value = "Apparemment, Tomcat utilise un Log4j2.xml. Pourquoi ? (pas un simple Log4j.properties),`nPourrais-je voir un exemple de Log4j2.xml bien implémenté ?"
$ws.Rows(11).RowHeight = 42.75

$ws.Range("B12").Value = "Test d'un controller Servlet"
$ws.Range("C12").Value = "Comment implémenter une org.apache….Request (héritant de HttpServletRequest) pour tester une méthode doGet(….) d'un controller Servlet ?"
$ws.Range("D12").Value = "Il pourrait être intéressant de tester un controller Servlet sans lancer l'appli web et donc Tomcat. Or, le conteneur de Servlet Tomcat fournit les HttpServletRequest et HttpServletResponse à passer aux méthodes doHead, doGet et doPost du controller Servlet. Comment instancier une HttpServletRequest dans un test unitaire JUnit en lui passant des paramètres ?"
$ws.Rows(12).RowHeight = 63.75

# ------------------------------------------------------------------
# 4) Finally, fill in the answer for the new "Choix du HTML" row.
# ------------------------------------------------------------------
$ws.Range("D5").Value = "ACAI préconise l'utilisation de "

# ------------------------------------------------------------------
# 5) Note: the stray path text that used to sit in (old) row 22 / D22
#    is already carried down to the new row 23 / D23 by the single
#    row-5 insert above (every row from the old row 10 onward shifts
#    down by exactly one), so no further action is required there.
# ------------------------------------------------------------------

# ------------------------------------------------------------------
# 6) Restore the active-cell selection shown in the file (E5).
# ------------------------------------------------------------------
$ws.Range("E5").Select()
